$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Increase Fractionated Coconut Oil (row 4) amount so the scaled mass becomes ~6 g,
# compensating by reducing the Water (row 2) amount so the overall blend total holds.
$ws.Range("B2").Value = 10.464957199674283
$ws.Range("B4").Value = 2.2102568019543289

# Apply a "0.0" number format to the Raw Amt column for the ingredient rows,
# matching the visual change made alongside the recipe update.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
$ws.Range("B5").NumberFormat = "0.0"

$ws.Range("H6").Select()
